$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")
$ws.Range("A2").Value = "Activity Test External Contact"
$ws.Range("B2").Value = "ActivityCompany"
$ws.Select()
$ws.Range("A2:B2").Select()
